$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# Update defined name ranges to reflect the new row extents for columns
# that gained (or need to reflect) additional command entries.
foreach ($n in $wb.Names) {
    switch ($n.Name) {
        "excel"     { $n.RefersTo = "='#system'!`$I`$2:`$I`$15" }
        "pdf"       { $n.RefersTo = "='#system'!`$S`$2:`$S`$17" }
        "web"       { $n.RefersTo = "='#system'!`$Z`$2:`$Z`$144" }
        "webcookie" { $n.RefersTo = "='#system'!`$AB`$2:`$AB`$10" }
    }
}


# Column I (excel)
$ws.Range("I2").Value = "assertPassword(file)"
$ws.Range("I3").Value = "clear(file,worksheet,range)"
$ws.Range("I4").Value = "clearPassword(file,password)"
$ws.Range("I5").Value = "columnarCsv(file,worksheet,ranges,output)"
$ws.Range("I6").Value = "csv(file,worksheet,range,output)"
$ws.Range("I7").Value = "json(file,worksheet,range,header,output)"
$ws.Range("I8").Value = "saveData(var,file,worksheet,range)"
$ws.Range("I9").Value = "saveRange(var,file,worksheet,range)"
$ws.Range("I10").Value = "saveTotalDataCount(file,worksheet,saveVar)"
$ws.Range("I11").Value = "setPassword(file,password)"
$ws.Range("I12").Value = "write(file,worksheet,startCell,data)"
$ws.Range("I13").Value = "writeAcross(file,worksheet,startCell,array)"
$ws.Range("I14").Value = "writeDown(file,worksheet,startCell,array)"
$ws.Range("I15").Value = "writeVar(var,file,worksheet,startCell)"

# Column S (pdf)
$ws.Range("S2").Value = "assertContentEqual(actualPdf,expectedPdf)"
$ws.Range("S3").Value = "assertFormElementPresent(var,name)"
$ws.Range("S4").Value = "assertFormValue(var,name,expected)"
$ws.Range("S5").Value = "assertFormValues(var,name,expectedValues,exactOrder)"
$ws.Range("S6").Value = "assertPatternNotPresent(pdf,regex)"
$ws.Range("S7").Value = "assertPatternPresent(pdf,regex)"
$ws.Range("S8").Value = "assertTextArray(pdf,textArray,ordered)"
$ws.Range("S9").Value = "assertTextNotPresent(pdf,text)"
$ws.Range("S10").Value = "assertTextPresent(pdf,text)"
$ws.Range("S11").Value = "count(pdf,text,var)"
$ws.Range("S12").Value = "saveAsPages(pdf,destination)"
$ws.Range("S13").Value = "saveAsPdf(profile,content,file)"
$ws.Range("S14").Value = "saveAsText(pdf,destination)"
$ws.Range("S15").Value = "saveFormValues(pdf,var,pageAndLineStartEnd,strategy)"
$ws.Range("S16").Value = "saveMetadata(pdf,var)"
$ws.Range("S17").Value = "saveToVar(pdf,var)"

# Column Z (web)
$ws.Range("Z2").Value = "assertAndClick(locator,label)"
$ws.Range("Z3").Value = "assertAttribute(locator,attrName,value)"
$ws.Range("Z4").Value = "assertAttributeContain(locator,attrName,contains)"
$ws.Range("Z5").Value = "assertAttributeNotContain(locator,attrName,contains)"
$ws.Range("Z6").Value = "assertAttributeNotPresent(locator,attrName)"
$ws.Range("Z7").Value = "assertAttributePresent(locator,attrName)"
$ws.Range("Z8").Value = "assertChecked(locator)"
$ws.Range("Z9").Value = "assertContainCount(locator,text,count)"
$ws.Range("Z10").Value = "assertCssNotPresent(locator,property)"
$ws.Range("Z11").Value = "assertCssPresent(locator,property,value)"
$ws.Range("Z12").Value = "assertElementByAttributes(nameValues)"
$ws.Range("Z13").Value = "assertElementByText(locator,text)"
$ws.Range("Z14").Value = "assertElementCount(locator,count)"
$ws.Range("Z15").Value = "assertElementEnabled(locator)"
$ws.Range("Z16").Value = "assertElementNotPresent(locator)"
$ws.Range("Z17").Value = "assertElementPresent(locator)"
$ws.Range("Z18").Value = "assertElementsPresent(prefix)"
$ws.Range("Z19").Value = "assertFocus(locator)"
$ws.Range("Z20").Value = "assertFrameCount(count)"
$ws.Range("Z21").Value = "assertFramePresent(frameName)"
$ws.Range("Z22").Value = "assertIECompatMode()"
$ws.Range("Z23").Value = "assertIENativeMode()"
$ws.Range("Z24").Value = "assertLinkByLabel(label)"
$ws.Range("Z25").Value = "assertMultiSelect(locator)"
$ws.Range("Z26").Value = "assertNotChecked(locator)"
$ws.Range("Z27").Value = "assertNotFocus(locator)"
$ws.Range("Z28").Value = "assertNotText(locator,text)"
$ws.Range("Z29").Value = "assertNotVisible(locator)"
$ws.Range("Z30").Value = "assertOneMatch(locator)"
$ws.Range("Z31").Value = "assertScrollbarHNotPresent(locator)"
$ws.Range("Z32").Value = "assertScrollbarHPresent(locator)"
$ws.Range("Z33").Value = "assertScrollbarVNotPresent(locator)"
$ws.Range("Z34").Value = "assertScrollbarVPresent(locator)"
$ws.Range("Z35").Value = "assertSingleSelect(locator)"
$ws.Range("Z36").Value = "assertTable(locator,row,column,text)"
$ws.Range("Z37").Value = "assertText(locator,text)"
$ws.Range("Z38").Value = "assertTextContains(locator,text)"
$ws.Range("Z39").Value = "assertTextCount(locator,text,count)"
$ws.Range("Z40").Value = "assertTextList(locator,list,ignoreOrder)"
$ws.Range("Z41").Value = "assertTextMatches(text,minMatch,scrollTo)"
$ws.Range("Z42").Value = "assertTextNotContain(locator,text)"
$ws.Range("Z43").Value = "assertTextNotPresent(text)"
$ws.Range("Z44").Value = "assertTextOrder(locator,descending)"
$ws.Range("Z45").Value = "assertTextPresent(text)"
$ws.Range("Z46").Value = "assertTitle(text)"
$ws.Range("Z47").Value = "assertValue(locator,value)"
$ws.Range("Z48").Value = "assertValueOrder(locator,descending)"
$ws.Range("Z49").Value = "assertVisible(locator)"
$ws.Range("Z50").Value = "checkAll(locator)"
$ws.Range("Z51").Value = "clearLocalStorage()"
$ws.Range("Z52").Value = "click(locator)"
$ws.Range("Z53").Value = "clickAll(locator)"
$ws.Range("Z54").Value = "clickAndWait(locator,waitMs)"
$ws.Range("Z55").Value = "clickByLabel(label)"
$ws.Range("Z56").Value = "clickByLabelAndWait(label,waitMs)"
$ws.Range("Z57").Value = "clickOffset(locator,x,y)"
$ws.Range("Z58").Value = "clickWithKeys(locator,keys)"
$ws.Range("Z59").Value = "close()"
$ws.Range("Z60").Value = "closeAll()"
$ws.Range("Z61").Value = "deselect(locator,text)"
$ws.Range("Z62").Value = "deselectMulti(locator,array)"
$ws.Range("Z63").Value = "dismissInvalidCert()"
$ws.Range("Z64").Value = "dismissInvalidCertPopup()"
$ws.Range("Z65").Value = "doubleClick(locator)"
$ws.Range("Z66").Value = "doubleClickAndWait(locator,waitMs)"
$ws.Range("Z67").Value = "doubleClickByLabel(label)"
$ws.Range("Z68").Value = "doubleClickByLabelAndWait(label,waitMs)"
$ws.Range("Z69").Value = "dragAndDrop(fromLocator,toLocator)"
$ws.Range("Z70").Value = "dragTo(fromLocator,xOffset,yOffset)"
$ws.Range("Z71").Value = "editLocalStorage(key,value)"
$ws.Range("Z72").Value = "executeScript(var,script)"
$ws.Range("Z73").Value = "focus(locator)"
$ws.Range("Z74").Value = "goBack()"
$ws.Range("Z75").Value = "goBackAndWait()"
$ws.Range("Z76").Value = "maximizeWindow()"
$ws.Range("Z77").Value = "mouseOver(locator)"
$ws.Range("Z78").Value = "open(url)"
$ws.Range("Z79").Value = "openAndWait(url,waitMs)"
$ws.Range("Z80").Value = "openHttpBasic(url,username,password)"
$ws.Range("Z81").Value = "openIgnoreTimeout(url)"
$ws.Range("Z82").Value = "refresh()"
$ws.Range("Z83").Value = "refreshAndWait()"
$ws.Range("Z84").Value = "resizeWindow(width,height)"
$ws.Range("Z85").Value = "rightClick(locator)"
$ws.Range("Z86").Value = "saveAllWindowIds(var)"
$ws.Range("Z87").Value = "saveAllWindowNames(var)"
$ws.Range("Z88").Value = "saveAttribute(var,locator,attrName)"
$ws.Range("Z89").Value = "saveAttributeList(var,locator,attrName)"
$ws.Range("Z90").Value = "saveBrowserVersion(var)"
$ws.Range("Z91").Value = "saveCount(var,locator)"
$ws.Range("Z92").Value = "saveDivsAsCsv(headers,rows,cells,nextPage,file)"
$ws.Range("Z93").Value = "saveElement(var,locator)"
$ws.Range("Z94").Value = "saveElements(var,locator)"
$ws.Range("Z95").Value = "saveInfiniteDivsAsCsv(config,file)"
$ws.Range("Z96").Value = "saveInfiniteTableAsCsv(config,file)"
$ws.Range("Z97").Value = "saveLocalStorage(var,key)"
$ws.Range("Z98").Value = "saveLocation(var)"
$ws.Range("Z99").Value = "savePageAs(var,sessionIdName,url)"
$ws.Range("Z100").Value = "savePageAsFile(sessionIdName,url,file)"
$ws.Range("Z101").Value = "saveSelectedText(var,locator)"
$ws.Range("Z102").Value = "saveSelectedValue(var,locator)"
$ws.Range("Z103").Value = "saveTableAsCsv(locator,nextPageLocator,file)"
$ws.Range("Z104").Value = "saveText(var,locator)"
$ws.Range("Z105").Value = "saveTextArray(var,locator)"
$ws.Range("Z106").Value = "saveTextSubstringAfter(var,locator,delim)"
$ws.Range("Z107").Value = "saveTextSubstringBefore(var,locator,delim)"
$ws.Range("Z108").Value = "saveTextSubstringBetween(var,locator,start,end)"
$ws.Range("Z109").Value = "saveTitle(var)"
$ws.Range("Z110").Value = "saveValue(var,locator)"
$ws.Range("Z111").Value = "saveValues(var,locator)"
$ws.Range("Z112").Value = "screenshot(file,locator)"
$ws.Range("Z113").Value = "scrollElement(locator,xOffset,yOffset)"
$ws.Range("Z114").Value = "scrollLeft(locator,pixel)"
$ws.Range("Z115").Value = "scrollPage(xOffset,yOffset)"
$ws.Range("Z116").Value = "scrollRight(locator,pixel)"
$ws.Range("Z117").Value = "scrollTo(locator)"
$ws.Range("Z118").Value = "select(locator,text)"
$ws.Range("Z119").Value = "selectAllOptions(locator)"
$ws.Range("Z120").Value = "selectFrame(locator)"
$ws.Range("Z121").Value = "selectMulti(locator,array)"
$ws.Range("Z122").Value = "selectMultiByValue(locator,array)"
$ws.Range("Z123").Value = "selectMultiOptions(locator)"
$ws.Range("Z124").Value = "selectText(locator)"
$ws.Range("Z125").Value = "selectWindow(winId)"
$ws.Range("Z126").Value = "selectWindowAndWait(winId,waitMs)"
$ws.Range("Z127").Value = "selectWindowByIndex(index)"
$ws.Range("Z128").Value = "selectWindowByIndexAndWait(index,waitMs)"
$ws.Range("Z129").Value = "switchBrowser(profile,config)"
$ws.Range("Z130").Value = "toggleSelections(locator)"
$ws.Range("Z131").Value = "type(locator,value)"
$ws.Range("Z132").Value = "typeKeys(locator,value)"
$ws.Range("Z133").Value = "uncheckAll(locator)"
$ws.Range("Z134").Value = "unselectAllText()"
$ws.Range("Z135").Value = "updateAttribute(locator,attrName,value)"
$ws.Range("Z136").Value = "upload(fieldLocator,file)"
$ws.Range("Z137").Value = "verifyContainText(locator,text)"
$ws.Range("Z138").Value = "verifyText(locator,text)"
$ws.Range("Z139").Value = "wait(waitMs)"
$ws.Range("Z140").Value = "waitForElementPresent(locator)"
$ws.Range("Z141").Value = "waitForElementsPresent(locators)"
$ws.Range("Z142").Value = "waitForPopUp(winId,waitMs)"
$ws.Range("Z143").Value = "waitForTextPresent(text)"
$ws.Range("Z144").Value = "waitForTitle(text)"

# Column AB (webcookie)
$ws.Range("AB2").Value = "assertNotPresent(name)"
$ws.Range("AB3").Value = "assertPresent(name)"
$ws.Range("AB4").Value = "assertValue(name,value)"
$ws.Range("AB5").Value = "clearCookieFields(var,remove)"
$ws.Range("AB6").Value = "delete(name)"
$ws.Range("AB7").Value = "deleteAll()"
$ws.Range("AB8").Value = "save(var,name)"
$ws.Range("AB9").Value = "saveAll(var)"
$ws.Range("AB10").Value = "saveAllAsText(var,exclude)"
